$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$origStyle = $cell.Style
$cell.Value = "'245.18"
$cell.Style = $origStyle

$cell = $ws.Range('D4')
$origStyle = $cell.Style
$cell.Value = "'5.412"
$cell.Style = $origStyle

$cell = $ws.Range('D6')
$origStyle = $cell.Style
$cell.Value = "'3.392"
$cell.Style = $origStyle

$cell = $ws.Range('D7')
$origStyle = $cell.Style
$cell.Value = "'0.8088"
$cell.Style = $origStyle

$cell = $ws.Range('D8')
$origStyle = $cell.Style
$cell.Value = "'0.9297"
$cell.Style = $origStyle

$ws.Range('B9').Value = 'WazirX'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$cell = $ws.Range('D9')
$origStyle = $cell.Style
$cell.Value = "'0.1419"
$cell.Style = $origStyle
$ws.Range('E9').Value = '8WazirXWRX'

$ws.Range('B10').Value = 'MandalaExchangeToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$cell = $ws.Range('D10')
$origStyle = $cell.Style
$cell.Value = "'0.07433"
$cell.Style = $origStyle
$ws.Range('E10').Value = '9MandalaExchangeTokenMDX'

$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$cell = $ws.Range('D11')
$origStyle = $cell.Style
$cell.Value = "'0.03378"
$cell.Style = $origStyle
$ws.Range('E11').Value = '10LiechtensteinCryptoassetsExchangeLCX'

$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$cell = $ws.Range('D12')
$origStyle = $cell.Style
$cell.Value = "'0.03026"
$cell.Style = $origStyle
$ws.Range('E12').Value = '11BitrueCoinBTR'

$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$cell = $ws.Range('D13')
$origStyle = $cell.Style
$cell.Value = "'0.09351"
$cell.Style = $origStyle
$ws.Range('E13').Value = '12BitMartTokenBMX'

$ws.Range('B14').Value = 'MCDex'
$ws.Range('C14').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$cell = $ws.Range('D14')
$origStyle = $cell.Style
$cell.Value = "'3.934"
$cell.Style = $origStyle
$ws.Range('E14').Value = '13MCDexMCB'

$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$cell = $ws.Range('D15')
$origStyle = $cell.Style
$cell.Value = "'0.001597"
$cell.Style = $origStyle
$ws.Range('E15').Value = '14BitForexTokenBF'

$ws.Range('B16').Value = 'CoinExToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$cell = $ws.Range('D16')
$origStyle = $cell.Style
$cell.Value = "'0.04805"
$cell.Style = $origStyle
$ws.Range('E16').Value = '15CoinExTokenCET'

$ws.Range('B17').Value = 'One'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$cell = $ws.Range('D17')
$origStyle = $cell.Style
$cell.Value = "'0.0005942"
$cell.Style = $origStyle
$ws.Range('E17').Value = '16OneONE'

$cell = $ws.Range('D18')
$origStyle = $cell.Style
$cell.Value = "'0.005358"
$cell.Style = $origStyle

$cell = $ws.Range('D20')
$origStyle = $cell.Style
$cell.Value = "'0.0009836"
$cell.Style = $origStyle

$cell = $ws.Range('D22')
$origStyle = $cell.Style
$cell.Value = "'3.660"
$cell.Style = $origStyle

$cell = $ws.Range('D23')
$origStyle = $cell.Style
$cell.Value = "'6.451"
$cell.Style = $origStyle

$cell = $ws.Range('D26')
$origStyle = $cell.Style
$cell.Value = "'0.1314"
$cell.Style = $origStyle

$cell = $ws.Range('D40')
$origStyle = $cell.Style
$cell.Value = "'0.03952"
$cell.Style = $origStyle

$ws.Range('B41').Value = 'BKEXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$cell = $ws.Range('D41')
$origStyle = $cell.Style
$cell.Value = "'0.1075"
$cell.Style = $origStyle
$ws.Range('E41').Value = '40BKEXTokenBKK'

$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$cell = $ws.Range('D42')
$origStyle = $cell.Style
$cell.Value = "'0.002721"
$cell.Style = $origStyle
$ws.Range('E42').Value = '41CEJICEJI'

$ws.Range('B43').Value = 'KickToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$cell = $ws.Range('D43')
$origStyle = $cell.Style
$cell.Value = "'0.003029"
$cell.Style = $origStyle
$ws.Range('E43').Value = '42KickTokenKICK'

$cell = $ws.Range('D44')
$origStyle = $cell.Style
$cell.Value = "'0.006815"
$cell.Style = $origStyle

$cell = $ws.Range('D45')
$origStyle = $cell.Style
$cell.Value = "'0.00005211"
$cell.Style = $origStyle

$ws.Range('E48').Value = '47CoinbaseStockTokenCOINBestin24h'

$cell = $ws.Range('D49')
$origStyle = $cell.Style
$cell.Value = "'0.002027"
$cell.Style = $origStyle
